# Update "Pais" sheet with refreshed COVID-19 country/provincia figures
# (data pulled at 11:09 instead of 09:52) and the resulting re-sort by
# "Casos totales" that reshuffled a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 28 de Agosto de 2020 a las 11:09"

function Set-Row($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

Set-Row 25  "Filipinas"   209544 3999 134474 71745 0 91  3325
Set-Row 26  "Indonesia"   165887 3003 120900 37818 0 105 7169

# Israel overtakes Ecuador
Set-Row 31  "Israel"      111493 1090 90158  20444 0 7   891
Set-Row 32  "Ecuador"     111219 0    95202  9546  0 0   6471

Set-Row 47  "Polonia"     65480  791  44785  18677 0 8   2018
Set-Row 52  "Singapur"    56666  94   55139  1500  0 0   27
Set-Row 54  "Barein"      50756  0    47370  3197  0 1   189
Set-Row 63  "Afganistan"  38140  11   29059  7679  0 1   1402
Set-Row 71  "Austria"     26590  229  22594  3263  0 0   733
Set-Row 73  "El Salvador" 25415  131  13570  11143 0 8   702
Set-Row 93  "Malasia"     9306   10   9030   151   0 0   125
Set-Row 111 "Hong Kong"   4769   13   4249   436   0 3   84

# Eslovaquia jumps ahead of Cabo Verde / Ruanda / Mozambique
Set-Row 119 "Eslovaquia"  3728   102  2225   1470  0 0   33
Set-Row 120 "Cabo Verde"  3699   0    2749   912   0 0   38
Set-Row 121 "Ruanda"      3672   0    1863   1794  0 0   15
Set-Row 122 "Mozambique"  3651   0    1968   1662  0 0   21

Set-Row 129 "Eslovenia"   2797   42   2236   428   0 0   133
Set-Row 154 "Letonia"     1375   9    1163   178   0 0   34

Write-Output "edit complete"
